$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the existing instance rows (2-4) down to (3-5) by copying values only,
# working bottom-up so data is not overwritten before it is copied, and the
# instance index in column A is renumbered to keep incrementing from 0.
for ($r = 4; $r -ge 2; $r--) {
    $dest = $r + 1
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
    $ws.Cells.Item($dest, 1).Value = $r - 1
}

# Populate the newly inserted row 2 with the new case-study instance data.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 10
$ws.Cells.Item(2, 3).Value = 18
$ws.Cells.Item(2, 4).Value = $false
$ws.Cells.Item(2, 5).Value = $false
$ws.Cells.Item(2, 6).Value = $false
$ws.Cells.Item(2, 7).Value = 9
$ws.Cells.Item(2, 8).Value = 0.34
$ws.Cells.Item(2, 9).Value = 3601.84
$ws.Cells.Item(2, 10).Value = 10.96

# Ensure the new last row's index cell (A5) carries the same formatting as the
# other index cells in column A (copied from row 4, which already has it).
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
